# changed codice_eyetr_museo to match participant folder names in all_gaze
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "2021_2bm"   = "2021_02bm"
    "2021_1bm"   = "2021_01bm"
    "2021_4bm"   = "2021_04bm"
    "2021_5bm"   = "2021_05bmnew"
    "2021_3bm"   = "2021_03bm"
    "2021_1bmf"  = "2021_01bmf"
    "2021_2bmf"  = "2021_02bmf"
    "2022_01bm"  = "2022_01bmnew"
    "2022_30bm"  = "2022_30bm_"
    "2022_41bm"  = "2022-41bm"
    "2022_38bm"  = "2022_38bmnew"
}

$usedRange = $ws.UsedRange
$rows = $usedRange.Rows.Count
$cols = $usedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $replacements.ContainsKey([string]$val)) {
            $cell.Value = $replacements[[string]$val]
        }
    }
}
